# Sacramento roster: two pairs of rows had their player data swapped while
# keeping the "No." index column (A) sequential:
#   - row 10 <-> row 11  (Terence Davis  <-> Chimezie Metu)
#   - row 17 <-> row 18  (Neemias Queta (TW) <-> Kessler Edwards)
# Column A (No.) is left untouched; columns B..K swap their contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-CellValue($ws, $addr) {
    return $ws.Range($addr).Value()
}

# Column I stores "years of experience" as text even though some values
# ("1","3","4",...) look numeric. A plain Range.Value assignment of such a
# string gets auto-coerced to a number by the General number format, so we
# briefly force a Text format, assign the value, then restore the original
# (default/"Normal") style so the cell's appearance/style is unaffected.
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

function Swap-RowContent($ws, $row1, $row2) {
    $plainCols = @("B", "C", "D", "E", "F", "G", "H", "J", "K")
    $textCols  = @("I")

    foreach ($col in $plainCols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = Get-CellValue $ws $addr1
        $v2 = Get-CellValue $ws $addr2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }

    foreach ($col in $textCols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = Get-CellValue $ws $addr1
        $v2 = Get-CellValue $ws $addr2
        Set-TextValue $ws $addr1 "$v2"
        Set-TextValue $ws $addr2 "$v1"
    }
}

Swap-RowContent $ws 10 11
Swap-RowContent $ws 17 18
